$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 83 - this shifts rows 83:103 down to 84:104,
# preserving all of their existing data/styles.
$ws.Rows.Item(83).Insert()

# Populate the newly inserted row 83 with the new record.
# Columns A,B,C,E,F,G,H,I,Q,R carry forward the same values as the rest of
# this data block (Femacal de La Calera / Coquimbo / Poroto granado row).
$ws.Cells.Item(83, 1).Value = 3
$ws.Cells.Item(83, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(83, 3).Value = "Coquimbo"
$ws.Cells.Item(83, 4).Value = 44508
$ws.Cells.Item(83, 5).Value = 5
$ws.Cells.Item(83, 6).Value = 100112030
$ws.Cells.Item(83, 7).Value = "Poroto granado"
$ws.Cells.Item(83, 8).Value = "Sin especificar"
$ws.Cells.Item(83, 9).Value = "Primera"
$ws.Cells.Item(83, 10).Value = 85
$ws.Cells.Item(83, 11).Value = 37000
$ws.Cells.Item(83, 12).Value = 38000
$ws.Cells.Item(83, 13).Value = 37529
$ws.Cells.Item(83, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(83, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(83, 16).Value = 1501
$ws.Cells.Item(83, 17).Value = 25
$ws.Cells.Item(83, 18).Value = "Hortaliza"
